# Rows 43/44, 54/55, 60/61 swap pairwise; rows 56/57/58 rotate as a 3-cycle
# (56 <- 57, 57 <- 58, 58 <- 56). Columns Y and AA ("Startdatum"/"Slutdatum")
# are left untouched because they hold identical text in every row touched
# here, so skipping them avoids the COM layer re-typing that text as a date.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Snapshot every affected row first, so no row is overwritten before it is read.
$save0_44 = $ws.Range("A44:X44").Value2
$save1_44 = $ws.Range("Z44:Z44").Value2
$save2_44 = $ws.Range("AB44:AY44").Value2
$save0_43 = $ws.Range("A43:X43").Value2
$save1_43 = $ws.Range("Z43:Z43").Value2
$save2_43 = $ws.Range("AB43:AY43").Value2
$save0_55 = $ws.Range("A55:X55").Value2
$save1_55 = $ws.Range("Z55:Z55").Value2
$save2_55 = $ws.Range("AB55:AY55").Value2
$save0_54 = $ws.Range("A54:X54").Value2
$save1_54 = $ws.Range("Z54:Z54").Value2
$save2_54 = $ws.Range("AB54:AY54").Value2
$save0_57 = $ws.Range("A57:X57").Value2
$save1_57 = $ws.Range("Z57:Z57").Value2
$save2_57 = $ws.Range("AB57:AY57").Value2
$save0_58 = $ws.Range("A58:X58").Value2
$save1_58 = $ws.Range("Z58:Z58").Value2
$save2_58 = $ws.Range("AB58:AY58").Value2
$save0_56 = $ws.Range("A56:X56").Value2
$save1_56 = $ws.Range("Z56:Z56").Value2
$save2_56 = $ws.Range("AB56:AY56").Value2
$save0_61 = $ws.Range("A61:X61").Value2
$save1_61 = $ws.Range("Z61:Z61").Value2
$save2_61 = $ws.Range("AB61:AY61").Value2
$save0_60 = $ws.Range("A60:X60").Value2
$save1_60 = $ws.Range("Z60:Z60").Value2
$save2_60 = $ws.Range("AB60:AY60").Value2

# 2) Write each destination row from its mapped source snapshot.
$ws.Range("A43:X43").Value2 = $save0_44
$ws.Range("Z43:Z43").Value2 = $save1_44
$ws.Range("AB43:AY43").Value2 = $save2_44
$ws.Range("A44:X44").Value2 = $save0_43
$ws.Range("Z44:Z44").Value2 = $save1_43
$ws.Range("AB44:AY44").Value2 = $save2_43
$ws.Range("A54:X54").Value2 = $save0_55
$ws.Range("Z54:Z54").Value2 = $save1_55
$ws.Range("AB54:AY54").Value2 = $save2_55
$ws.Range("A55:X55").Value2 = $save0_54
$ws.Range("Z55:Z55").Value2 = $save1_54
$ws.Range("AB55:AY55").Value2 = $save2_54
$ws.Range("A56:X56").Value2 = $save0_57
$ws.Range("Z56:Z56").Value2 = $save1_57
$ws.Range("AB56:AY56").Value2 = $save2_57
$ws.Range("A57:X57").Value2 = $save0_58
$ws.Range("Z57:Z57").Value2 = $save1_58
$ws.Range("AB57:AY57").Value2 = $save2_58
$ws.Range("A58:X58").Value2 = $save0_56
$ws.Range("Z58:Z58").Value2 = $save1_56
$ws.Range("AB58:AY58").Value2 = $save2_56
$ws.Range("A60:X60").Value2 = $save0_61
$ws.Range("Z60:Z60").Value2 = $save1_61
$ws.Range("AB60:AY60").Value2 = $save2_61
$ws.Range("A61:X61").Value2 = $save0_60
$ws.Range("Z61:Z61").Value2 = $save1_60
$ws.Range("AB61:AY61").Value2 = $save2_60
